$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 28574158
$ws.Range("I100").Value = 1775
$ws.Range("J100").Value = 66670668
$ws.Range("K100").Value = 1775
$ws.Range("L100").Value = 66670668
$ws.Range("M100").Value = -1234
$ws.Range("N100").Value = -66671750

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2716739.8
$ws.Range("I61").Value = 1437603
$ws.Range("J61").Value = 7353611
$ws.Range("K61").Value = 1437603
$ws.Range("L61").Value = 7353611
$ws.Range("M61").Value = -1437391
$ws.Range("N61").Value = -7354035
$ws.Range("H110").Value = 1055.2858
$ws.Range("I110").Value = 1055.2858
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1055.2858
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 989.7141999999999
$ws.Range("N110").Value = $null
$ws.Range("H112").Value = 34346
$ws.Range("J112").Value = 34346
$ws.Range("L112").Value = 34346
$ws.Range("N112").Value = -37300
$ws.Range("H122").Value = 2318
$ws.Range("I122").Value = 1146.6
$ws.Range("K122").Value = 3439.8
$ws.Range("M122").Value = -989.7999999999997
$ws.Range("H123").Value = 53677.8
$ws.Range("J123").Value = 53677.8
$ws.Range("L123").Value = 53677.8
$ws.Range("N123").Value = -63477.8
$ws.Range("H132").Value = 12549181
$ws.Range("I132").Value = 12350291
$ws.Range("J132").Value = 13891688
$ws.Range("K132").Value = 37050873
$ws.Range("L132").Value = 41675064
$ws.Range("M132").Value = -37048343
$ws.Range("N132").Value = -41680124
$ws.Range("H136").Value = 2716739.8
$ws.Range("I136").Value = 1437603
$ws.Range("J136").Value = 7353611
$ws.Range("K136").Value = 4312809
$ws.Range("L136").Value = 22060833
$ws.Range("M136").Value = -4310259
$ws.Range("N136").Value = -22065933

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1771.4286
$ws.Range("I105").Value = 1466.6666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1466.6666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 280.3334
$ws.Range("N105").Value = -5494
$ws.Range("H110").Value = 27400
$ws.Range("J110").Value = 27400
$ws.Range("L110").Value = 27400
$ws.Range("N110").Value = -35580
$ws.Range("H134").Value = 7339832
$ws.Range("I134").Value = 8197771
$ws.Range("J134").Value = 2978642.8
$ws.Range("K134").Value = 24593313
$ws.Range("L134").Value = 8935928.399999999
$ws.Range("M134").Value = -24590778
$ws.Range("N134").Value = -8940998.399999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3364437.8
$ws.Range("I31").Value = 2316805.8
$ws.Range("J31").Value = 4815005
$ws.Range("K31").Value = 2316805.8
$ws.Range("L31").Value = 4815005
$ws.Range("M31").Value = -2316510.8
$ws.Range("N31").Value = -4815595
$ws.Range("H34").Value = 3364437.8
$ws.Range("I34").Value = 2316805.8
$ws.Range("J34").Value = 4815005
$ws.Range("K34").Value = 2316805.8
$ws.Range("L34").Value = 4815005
$ws.Range("M34").Value = -2316603.8
$ws.Range("N34").Value = -4815409
$ws.Range("H58").Value = 4062606.8
$ws.Range("I58").Value = 2236880.8
$ws.Range("J58").Value = 11365511
$ws.Range("K58").Value = 2236880.8
$ws.Range("L58").Value = 11365511
$ws.Range("M58").Value = -2236677.8
$ws.Range("N58").Value = -11365917
$ws.Range("H132").Value = 2780542.8
$ws.Range("I132").Value = 4547448.5
$ws.Range("J132").Value = 3976
$ws.Range("K132").Value = 13642345.5
$ws.Range("L132").Value = 11928
$ws.Range("M132").Value = -13639815.5
$ws.Range("N132").Value = -16988
$ws.Range("H134").Value = 1147810.1
$ws.Range("I134").Value = 4965.8076
$ws.Range("J134").Value = 4449360
$ws.Range("K134").Value = 14897.4228
$ws.Range("L134").Value = 13348080
$ws.Range("M134").Value = -12362.4228
$ws.Range("N134").Value = -13353150
$ws.Range("H136").Value = 4062606.8
$ws.Range("I136").Value = 2236880.8
$ws.Range("J136").Value = 11365511
$ws.Range("K136").Value = 6710642.399999999
$ws.Range("L136").Value = 34096533
$ws.Range("M136").Value = -6708092.399999999
$ws.Range("N136").Value = -34101633

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1314.625
$ws.Range("I2").Value = 83.333336
$ws.Range("J2").Value = 2053.4
$ws.Range("K2").Value = 500.000016
$ws.Range("L2").Value = 12320.4
$ws.Range("M2").Value = -387.000016
$ws.Range("N2").Value = -12546.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12824274
$ws.Range("I122").Value = 4080
$ws.Range("J122").Value = 83335340
$ws.Range("K122").Value = 12240
$ws.Range("L122").Value = 250006020
$ws.Range("M122").Value = -9790
$ws.Range("N122").Value = -250010920

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4584.421
$ws.Range("I22").Value = 3527.1428
$ws.Range("J22").Value = 5201.1665
$ws.Range("K22").Value = 3527.1428
$ws.Range("L22").Value = 5201.1665
$ws.Range("M22").Value = -3232.1428
$ws.Range("N22").Value = -5791.1665
$ws.Range("H27").Value = 4584.421
$ws.Range("I27").Value = 3527.1428
$ws.Range("J27").Value = 5201.1665
$ws.Range("K27").Value = 3527.1428
$ws.Range("L27").Value = 5201.1665
$ws.Range("M27").Value = -3420.1428
$ws.Range("N27").Value = -5415.1665
$ws.Range("H40").Value = 1828.4615
$ws.Range("I40").Value = 1138.7
$ws.Range("J40").Value = 2259.5625
$ws.Range("K40").Value = 1138.7
$ws.Range("L40").Value = 2259.5625
$ws.Range("M40").Value = -1002.7
$ws.Range("N40").Value = -2531.5625
$ws.Range("H110").Value = 28700
$ws.Range("J110").Value = 28700
$ws.Range("L110").Value = 28700
$ws.Range("N110").Value = -36880
$ws.Range("H136").Value = 1369171.1
$ws.Range("I136").Value = 1783590.5
$ws.Range("J136").Value = 1587.4
$ws.Range("K136").Value = 5350771.5
$ws.Range("L136").Value = 4762.200000000001
$ws.Range("M136").Value = -5348221.5
$ws.Range("N136").Value = -9862.200000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = $null
$ws.Range("N38").Value = $null
$ws.Range("H39").Value = 3900
$ws.Range("I39").Value = 3900
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3900
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -3487
$ws.Range("N39").Value = $null
$ws.Range("H47").Value = 17980
$ws.Range("J47").Value = 17980
$ws.Range("L47").Value = 17980
$ws.Range("N47").Value = -19124
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = $null
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("H132").Value = 1158963.4
$ws.Range("I132").Value = 822976.4399999999
$ws.Range("J132").Value = 2166924.2
$ws.Range("K132").Value = 2468929.32
$ws.Range("L132").Value = 6500772.600000001
$ws.Range("M132").Value = -2466399.32
$ws.Range("N132").Value = -6505832.600000001
